$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 8857.615  # H64: 8919.191999999999 -> 8857.615
$ws.Cells.Item(64, 10).Value = 9559.046  # J64: 9631.817999999999 -> 9559.046
$ws.Cells.Item(64, 12).Value = 9559.046  # L64: 9631.817999999999 -> 9559.046
$ws.Cells.Item(64, 14).Value = -10055.046  # N64: -10127.818 -> -10055.046
$ws.Cells.Item(67, 8).Value = 8857.615  # H67: 8919.191999999999 -> 8857.615
$ws.Cells.Item(67, 10).Value = 9559.046  # J67: 9631.817999999999 -> 9559.046
$ws.Cells.Item(67, 12).Value = 9559.046  # L67: 9631.817999999999 -> 9559.046
$ws.Cells.Item(67, 14).Value = -11275.046  # N67: -11347.818 -> -11275.046
$ws.Cells.Item(97, 8).Value = 5857.2  # H97: 7146.5 -> 5857.2
$ws.Cells.Item(97, 10).Value = 5857.2  # J97: 7146.5 -> 5857.2
$ws.Cells.Item(97, 12).Value = 17571.6  # L97: 21439.5 -> 17571.6
$ws.Cells.Item(97, 14).Value = -18563.6  # N97: -22431.5 -> -18563.6
$ws.Cells.Item(130, 8).Value = 75000  # H130: 0 -> 75000
$ws.Cells.Item(130, 10).Value = 75000  # J130: 0 -> 75000
$ws.Cells.Item(130, 12).Value = 75000  # L130: 0 -> 75000
$ws.Cells.Item(130, 14).Value = -85040  # N130: None -> -85040
$ws.Cells.Item(132, 8).Value = 2626.8708  # H132: 2763.2068 -> 2626.8708
$ws.Cells.Item(132, 9).Value = 1896.64  # I132: 2005.0435 -> 1896.64
$ws.Cells.Item(132, 11).Value = 5689.92  # K132: 6015.1305 -> 5689.92
$ws.Cells.Item(132, 13).Value = -3159.92  # M132: -3485.1305 -> -3159.92
$ws.Cells.Item(137, 8).Value = 7962.154  # H137: 8244.719999999999 -> 7962.154
$ws.Cells.Item(137, 9).Value = 11044.333  # I137: 12589.77 -> 11044.333
$ws.Cells.Item(137, 10).Value = 3759.182  # J137: 3537.5833 -> 3759.182
$ws.Cells.Item(137, 11).Value = 33132.999  # K137: 37769.31 -> 33132.999
$ws.Cells.Item(137, 12).Value = 11277.546  # L137: 10612.7499 -> 11277.546
$ws.Cells.Item(137, 13).Value = -30582.999  # M137: -35219.31 -> -30582.999
$ws.Cells.Item(137, 14).Value = -16377.546  # N137: -15712.7499 -> -16377.546
$ws.Cells.Item(141, 8).Value = 3499.25  # H141: 3611.875 -> 3499.25
$ws.Cells.Item(141, 9).Value = 3598.8  # I141: 4248.75 -> 3598.8
$ws.Cells.Item(141, 10).Value = 3333.3333  # J141: 2975 -> 3333.3333
$ws.Cells.Item(141, 11).Value = 10796.4  # K141: 12746.25 -> 10796.4
$ws.Cells.Item(141, 12).Value = 9999.999899999999  # L141: 8925 -> 9999.999899999999
$ws.Cells.Item(141, 13).Value = -5616.400000000001  # M141: -7566.25 -> -5616.400000000001
$ws.Cells.Item(141, 14).Value = -20359.9999  # N141: -19285 -> -20359.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 125  # H5: 107.333336 -> 125
$ws.Cells.Item(5, 10).Value = 0  # J5: 72 -> 0
$ws.Cells.Item(5, 12).Value = 0  # L5: 72 -> 0
$ws.Cells.Item(5, 14).ClearContents()  # N5: was -296
$ws.Cells.Item(32, 8).Value = 2072.4849  # H32: 2174.9707 -> 2072.4849
$ws.Cells.Item(32, 9).Value = 2060.1587  # I32: 2167.754 -> 2060.1587
$ws.Cells.Item(32, 11).Value = 2060.1587  # K32: 2167.754 -> 2060.1587
$ws.Cells.Item(32, 13).Value = -1773.1587  # M32: -1880.754 -> -1773.1587
$ws.Cells.Item(63, 8).Value = 5985  # H63: 7135.7144 -> 5985
$ws.Cells.Item(63, 9).Value = 2973.75  # I63: 3316.6667 -> 2973.75
$ws.Cells.Item(63, 11).Value = 2973.75  # K63: 3316.6667 -> 2973.75
$ws.Cells.Item(63, 13).Value = -2287.75  # M63: -2630.6667 -> -2287.75
$ws.Cells.Item(66, 8).Value = 5985  # H66: 7135.7144 -> 5985
$ws.Cells.Item(66, 9).Value = 2973.75  # I66: 3316.6667 -> 2973.75
$ws.Cells.Item(66, 11).Value = 14868.75  # K66: 16583.3335 -> 14868.75
$ws.Cells.Item(66, 13).Value = -11436.75  # M66: -13151.3335 -> -11436.75
$ws.Cells.Item(68, 8).Value = 0  # H68: 40000 -> 0
$ws.Cells.Item(68, 10).Value = 0  # J68: 40000 -> 0
$ws.Cells.Item(68, 12).Value = 0  # L68: 40000 -> 0
$ws.Cells.Item(68, 14).ClearContents()  # N68: was -41622
$ws.Cells.Item(71, 8).Value = 0  # H71: 40000 -> 0
$ws.Cells.Item(71, 10).Value = 0  # J71: 40000 -> 0
$ws.Cells.Item(71, 12).Value = 0  # L71: 120000 -> 0
$ws.Cells.Item(71, 14).ClearContents()  # N71: was -128112
$ws.Cells.Item(74, 8).Value = 2850.3  # H74: 2527.3914 -> 2850.3
$ws.Cells.Item(74, 9).Value = 2313.0625  # I74: 2007 -> 2313.0625
$ws.Cells.Item(74, 11).Value = 2313.0625  # K74: 2007 -> 2313.0625
$ws.Cells.Item(74, 13).Value = -1439.0625  # M74: -1133 -> -1439.0625
$ws.Cells.Item(77, 8).Value = 2850.3  # H77: 2527.3914 -> 2850.3
$ws.Cells.Item(77, 9).Value = 2313.0625  # I77: 2007 -> 2313.0625
$ws.Cells.Item(77, 11).Value = 11565.3125  # K77: 10035 -> 11565.3125
$ws.Cells.Item(77, 13).Value = -7197.3125  # M77: -5667 -> -7197.3125
$ws.Cells.Item(88, 8).Value = 3664.8462  # H88: 3488.7856 -> 3664.8462
$ws.Cells.Item(88, 9).Value = 5000  # I88: 4050 -> 5000
$ws.Cells.Item(88, 11).Value = 5000  # K88: 4050 -> 5000
$ws.Cells.Item(88, 13).Value = -4594  # M88: -3644 -> -4594
$ws.Cells.Item(91, 8).Value = 3664.8462  # H91: 3488.7856 -> 3664.8462
$ws.Cells.Item(91, 9).Value = 5000  # I91: 4050 -> 5000
$ws.Cells.Item(91, 11).Value = 5000  # K91: 4050 -> 5000
$ws.Cells.Item(91, 13).Value = -3596  # M91: -2646 -> -3596
$ws.Cells.Item(102, 8).Value = 5891.467  # H102: 3971.111 -> 5891.467
$ws.Cells.Item(102, 9).Value = 3152.4443  # I102: 2311 -> 3152.4443
$ws.Cells.Item(102, 10).Value = 10000  # J102: 8714.286 -> 10000
$ws.Cells.Item(102, 11).Value = 3152.4443  # K102: 2311 -> 3152.4443
$ws.Cells.Item(102, 12).Value = 10000  # L102: 8714.286 -> 10000
$ws.Cells.Item(102, 13).Value = -1530.4443  # M102: -689 -> -1530.4443
$ws.Cells.Item(102, 14).Value = -13244  # N102: -11958.286 -> -13244
$ws.Cells.Item(132, 8).Value = 2618.7317  # H132: 2533.1396 -> 2618.7317
$ws.Cells.Item(132, 9).Value = 2438.0833  # I132: 2310.1794 -> 2438.0833
$ws.Cells.Item(132, 10).Value = 3919.4  # J132: 4707 -> 3919.4
$ws.Cells.Item(132, 11).Value = 7314.249899999999  # K132: 6930.5382 -> 7314.249899999999
$ws.Cells.Item(132, 12).Value = 11758.2  # L132: 14121 -> 11758.2
$ws.Cells.Item(132, 13).Value = -4784.249899999999  # M132: -4400.5382 -> -4784.249899999999
$ws.Cells.Item(132, 14).Value = -16818.2  # N132: -19181 -> -16818.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 125  # H4: 107.333336 -> 125
$ws.Cells.Item(4, 10).Value = 0  # J4: 72 -> 0
$ws.Cells.Item(4, 12).Value = 0  # L4: 72 -> 0
$ws.Cells.Item(4, 14).ClearContents()  # N4: was -302

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 2275.1667  # H2: 2454.818 -> 2275.1667
$ws.Cells.Item(2, 10).Value = 1474.75  # J2: 1866.6666 -> 1474.75
$ws.Cells.Item(2, 12).Value = 1474.75  # L2: 1866.6666 -> 1474.75
$ws.Cells.Item(2, 14).Value = -1700.75  # N2: -2092.6666 -> -1700.75
$ws.Cells.Item(3, 8).Value = 0  # H3: 1500 -> 0
$ws.Cells.Item(3, 10).Value = 0  # J3: 1500 -> 0
$ws.Cells.Item(3, 12).Value = 0  # L3: 1500 -> 0
$ws.Cells.Item(3, 14).ClearContents()  # N3: was -1726
$ws.Cells.Item(4, 8).Value = 7375  # H4: 9750 -> 7375
$ws.Cells.Item(4, 10).Value = 6666.6665  # J4: 10000 -> 6666.6665
$ws.Cells.Item(4, 12).Value = 6666.6665  # L4: 10000 -> 6666.6665
$ws.Cells.Item(4, 14).Value = -6890.6665  # N4: -10224 -> -6890.6665
$ws.Cells.Item(10, 8).Value = 500  # H10: 0 -> 500
$ws.Cells.Item(10, 10).Value = 500  # J10: 0 -> 500
$ws.Cells.Item(10, 12).Value = 500  # L10: 0 -> 500
$ws.Cells.Item(10, 14).Value = -778  # N10: None -> -778
$ws.Cells.Item(22, 8).Value = 1444.409  # H22: 1340.2916 -> 1444.409
$ws.Cells.Item(22, 9).Value = 985.5333000000001  # I22: 936.4375 -> 985.5333000000001
$ws.Cells.Item(22, 10).Value = 2427.7144  # J22: 2148 -> 2427.7144
$ws.Cells.Item(22, 11).Value = 985.5333000000001  # K22: 936.4375 -> 985.5333000000001
$ws.Cells.Item(22, 12).Value = 2427.7144  # L22: 2148 -> 2427.7144
$ws.Cells.Item(22, 13).Value = -635.5333000000001  # M22: -586.4375 -> -635.5333000000001
$ws.Cells.Item(22, 14).Value = -3127.7144  # N22: -2848 -> -3127.7144
$ws.Cells.Item(132, 8).Value = 2674.0588  # H132: 2616.5293 -> 2674.0588
$ws.Cells.Item(132, 9).Value = 2634.1333  # I132: 2594.1875 -> 2634.1333
$ws.Cells.Item(132, 10).Value = 2973.5  # J132: 2974 -> 2973.5
$ws.Cells.Item(132, 11).Value = 7902.3999  # K132: 7782.5625 -> 7902.3999
$ws.Cells.Item(132, 12).Value = 8920.5  # L132: 8922 -> 8920.5
$ws.Cells.Item(132, 13).Value = -5372.3999  # M132: -5252.5625 -> -5372.3999
$ws.Cells.Item(132, 14).Value = -13980.5  # N132: -13982 -> -13980.5
$ws.Cells.Item(134, 8).Value = 2144.0386  # H134: 2174.84 -> 2144.0386
$ws.Cells.Item(134, 9).Value = 1955.3636  # I134: 1983.0476 -> 1955.3636
$ws.Cells.Item(134, 11).Value = 5866.0908  # K134: 5949.142800000001 -> 5866.0908
$ws.Cells.Item(134, 13).Value = -3331.0908  # M134: -3414.142800000001 -> -3331.0908
$ws.Cells.Item(140, 8).Value = 149999  # H140: 136662.67 -> 149999
$ws.Cells.Item(140, 10).Value = 149999  # J140: 136662.67 -> 149999
$ws.Cells.Item(140, 12).Value = 149999  # L140: 136662.67 -> 149999
$ws.Cells.Item(140, 14).Value = -160359  # N140: -147022.67 -> -160359

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2002.762  # H5: 2088.95 -> 2002.762
$ws.Cells.Item(5, 10).Value = 3496.4546  # J5: 3818.2 -> 3496.4546
$ws.Cells.Item(5, 12).Value = 10489.3638  # L5: 11454.6 -> 10489.3638
$ws.Cells.Item(5, 14).Value = -10713.3638  # N5: -11678.6 -> -10713.3638
$ws.Cells.Item(135, 8).Value = 2002.762  # H135: 2088.95 -> 2002.762
$ws.Cells.Item(135, 10).Value = 3496.4546  # J135: 3818.2 -> 3496.4546
$ws.Cells.Item(135, 12).Value = 31468.0914  # L135: 34363.8 -> 31468.0914
$ws.Cells.Item(135, 14).Value = -36538.0914  # N135: -39433.8 -> -36538.0914

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 10244.917  # H80: 12846.117 -> 10244.917
$ws.Cells.Item(80, 9).Value = 16567.666  # I80: 20230.572 -> 16567.666
$ws.Cells.Item(80, 10).Value = 6451.2666  # J80: 7677 -> 6451.2666
$ws.Cells.Item(80, 11).Value = 16567.666  # K80: 20230.572 -> 16567.666
$ws.Cells.Item(80, 12).Value = 6451.2666  # L80: 7677 -> 6451.2666
$ws.Cells.Item(80, 13).Value = -15569.666  # M80: -19232.572 -> -15569.666
$ws.Cells.Item(80, 14).Value = -8447.266599999999  # N80: -9673 -> -8447.266599999999
$ws.Cells.Item(83, 8).Value = 10244.917  # H83: 12846.117 -> 10244.917
$ws.Cells.Item(83, 9).Value = 16567.666  # I83: 20230.572 -> 16567.666
$ws.Cells.Item(83, 10).Value = 6451.2666  # J83: 7677 -> 6451.2666
$ws.Cells.Item(83, 11).Value = 82838.33  # K83: 101152.86 -> 82838.33
$ws.Cells.Item(83, 12).Value = 32256.333  # L83: 38385 -> 32256.333
$ws.Cells.Item(83, 13).Value = -77846.33  # M83: -96160.86 -> -77846.33
$ws.Cells.Item(83, 14).Value = -42240.333  # N83: -48369 -> -42240.333
$ws.Cells.Item(122, 8).Value = 6909.0454  # H122: 7432.8887 -> 6909.0454
$ws.Cells.Item(122, 9).Value = 7138.1875  # I122: 7731.077 -> 7138.1875
$ws.Cells.Item(122, 10).Value = 6298  # J122: 6657.6 -> 6298
$ws.Cells.Item(122, 11).Value = 21414.5625  # K122: 23193.231 -> 21414.5625
$ws.Cells.Item(122, 12).Value = 18894  # L122: 19972.8 -> 18894
$ws.Cells.Item(122, 13).Value = -18964.5625  # M122: -20743.231 -> -18964.5625
$ws.Cells.Item(122, 14).Value = -23794  # N122: -24872.8 -> -23794
$ws.Cells.Item(132, 8).Value = 3226.795  # H132: 3717.8276 -> 3226.795
$ws.Cells.Item(132, 9).Value = 3078.2  # I132: 3560.5454 -> 3078.2
$ws.Cells.Item(132, 10).Value = 3722.111  # J132: 4212.143 -> 3722.111
$ws.Cells.Item(132, 11).Value = 9234.599999999999  # K132: 10681.6362 -> 9234.599999999999
$ws.Cells.Item(132, 12).Value = 11166.333  # L132: 12636.429 -> 11166.333
$ws.Cells.Item(132, 13).Value = -6704.599999999999  # M132: -8151.636200000001 -> -6704.599999999999
$ws.Cells.Item(132, 14).Value = -16226.333  # N132: -17696.429 -> -16226.333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6568.143  # H7: 6664 -> 6568.143
$ws.Cells.Item(7, 9).Value = 6663.6665  # I7: 6797.8 -> 6663.6665
$ws.Cells.Item(7, 11).Value = 6663.6665  # K7: 6797.8 -> 6663.6665
$ws.Cells.Item(7, 13).Value = -6551.6665  # M7: -6685.8 -> -6551.6665
$ws.Cells.Item(46, 8).Value = 2734.2083  # H46: 2778.3635 -> 2734.2083
$ws.Cells.Item(46, 10).Value = 3157.2354  # J46: 3278.4 -> 3157.2354
$ws.Cells.Item(46, 12).Value = 3157.2354  # L46: 3278.4 -> 3157.2354
$ws.Cells.Item(46, 14).Value = -3533.2354  # N46: -3654.4 -> -3533.2354
$ws.Cells.Item(126, 8).Value = 6568.143  # H126: 6664 -> 6568.143
$ws.Cells.Item(126, 9).Value = 6663.6665  # I126: 6797.8 -> 6663.6665
$ws.Cells.Item(126, 11).Value = 19990.9995  # K126: 20393.4 -> 19990.9995
$ws.Cells.Item(126, 13).Value = -17520.9995  # M126: -17923.4 -> -17520.9995
$ws.Cells.Item(132, 8).Value = 3342  # H132: 3338.1538 -> 3342
$ws.Cells.Item(132, 9).Value = 3202.913  # I132: 3198.5652 -> 3202.913
$ws.Cells.Item(132, 11).Value = 9608.739  # K132: 9595.695599999999 -> 9608.739
$ws.Cells.Item(132, 13).Value = -7078.739  # M132: -7065.695599999999 -> -7078.739
$ws.Cells.Item(136, 8).Value = 6428.4287  # H136: 4996 -> 6428.4287
$ws.Cells.Item(136, 9).Value = 6428.4287  # I136: 4996 -> 6428.4287
$ws.Cells.Item(136, 11).Value = 19285.2861  # K136: 14988 -> 19285.2861
$ws.Cells.Item(136, 13).Value = -16735.2861  # M136: -12438 -> -16735.2861

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 2292.6843  # H100: 2308.4211 -> 2292.6843
$ws.Cells.Item(100, 9).Value = 1972.4375  # I100: 1991.125 -> 1972.4375
$ws.Cells.Item(100, 11).Value = 3944.875  # K100: 3982.25 -> 3944.875
$ws.Cells.Item(100, 13).Value = -3403.875  # M100: -3441.25 -> -3403.875
$ws.Cells.Item(132, 8).Value = 4958.609  # H132: 5106.9546 -> 4958.609
$ws.Cells.Item(132, 9).Value = 4309.75  # I132: 4484.067 -> 4309.75
$ws.Cells.Item(132, 11).Value = 12929.25  # K132: 13452.201 -> 12929.25
$ws.Cells.Item(132, 13).Value = -10399.25  # M132: -10922.201 -> -10399.25
$ws.Cells.Item(136, 8).Value = 1228.4166  # H136: 1289.3429 -> 1228.4166
$ws.Cells.Item(136, 9).Value = 1126.5807  # I136: 1178 -> 1126.5807
$ws.Cells.Item(136, 10).Value = 1859.8  # J136: 1827.5 -> 1859.8
$ws.Cells.Item(136, 11).Value = 3379.7421  # K136: 3534 -> 3379.7421
$ws.Cells.Item(136, 12).Value = 5579.4  # L136: 5482.5 -> 5579.4
$ws.Cells.Item(136, 13).Value = -829.7420999999999  # M136: -984 -> -829.7420999999999
$ws.Cells.Item(136, 14).Value = -10679.4  # N136: -10582.5 -> -10679.4
